# Generate Report for Handoff
#
# A fresh localization hand-off run just completed for the four files that
# were sitting in "Ready for handoff" (93665ec8…, b243efce…, fa018c9e…,
# fb3af559…). Their Priority flips from "low" to "ht", and the "Latest
# Handoff Datetime" for that batch is refreshed to the new run's timestamp,
# once per locale sheet. The Overview sheet's "Latest HO Xliff Generate
# Date" mirrors the de-de locale's handoff datetime, so it is refreshed too.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7 correspond to 93665ec8, b243efce, fa018c9e, fb3af559
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-08-12 22:38:07"

# de-de: same four rows
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-08-12 22:38:15"

# Overview sheet's "Latest HO Xliff Generate Date" column for the same
# four files mirrors the de-de handoff datetime.
$wsOverview.Range("G4:G7").Value = "2016-08-12 22:38:15"
